$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of kaspa buy data for 2025-10-24.
# Set the date cell as text first (so Excel doesn't silently convert the
# "MM/DD/YYYY" looking string into a date serial number + date format),
# then drop the style back to Normal so the cell carries no special
# formatting, matching the other plain date-string cells in this sheet.
$dateCell = $ws.Cells.Item(11, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "10/24/2025"
$dateCell.Style = "Normal"

$ws.Cells.Item(11, 2).Value = 482.9750000000004
$ws.Cells.Item(11, 3).Value = 0.1035250271753196
$ws.Cells.Item(11, 4).Value = 25
